$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ G=32.669943; H=98.009829; I=0.344160451822467; J=0.3441604518224671; M=68.18146900000001; N=204.544407; O=0.1244286043321187; P=0.1244286043321187; Q=2227.484705886267; R=20047.3623529764; S=0.04282340468658093; T=0.04282340468658094 }
    3  = @{ G=32.669943; H=98.009829; I=0.344160451822467; J=0.3441604518224671; O=0.345973452289334; P=0.3459734522893341; Q=6193.516175430058; R=55741.64557887053; S=0.1190703796584759; T=0.119070379658476 }
    4  = @{ G=32.669943; H=98.009829; I=0.344160451822467; J=0.3441604518224671; M=188.0130056666667; N=564.0390170000001; O=0.3431166302883566; P=0.3431166302883567; Q=6142.374178388677; R=55281.3676054981; S=0.1180871745078432; T=0.1180871745078432 }
    5  = @{ G=32.669943; H=98.009829; I=0.344160451822467; J=0.3441604518224671; M=102.1836573333333; N=306.550972; O=0.1864813130901906; P=0.1864813130901907; Q=3338.334260611532; R=30045.00834550379; S=0.06417949296956695; T=0.06417949296956696 }
    6  = @{ I=0.4225581378297699; J=0.4225581378297699; M=68.18146900000001; N=204.544407; O=0.1244286043321187; P=0.1244286043321187; Q=2734.892357269241; R=24614.03121542316; S=0.05257831933933731; T=0.05257831933933731 }
    7  = @{ I=0.4225581378297699; J=0.4225581378297699; O=0.345973452289334; P=0.3459734522893341; S=0.1461938977379177; T=0.1461938977379177 }
    8  = @{ I=0.4225581378297699; J=0.4225581378297699; M=188.0130056666667; N=564.0390170000001; O=0.3431166302883566; P=0.3431166302883567; Q=7541.570162781108; R=67874.13146502998; S=0.1449867243530736; T=0.1449867243530736 }
    9  = @{ I=0.4225581378297699; J=0.4225581378297699; M=102.1836573333333; N=306.550972; O=0.1864813130901906; P=0.1864813130901907; Q=4098.786775608373; R=36889.08098047536; S=0.07879919639944126; T=0.07879919639944127 }
    10 = @{ G=21.811182; H=65.43354600000001; I=0.2297691872894318; J=0.2297691872894318; M=68.18146900000001; N=204.544407; O=0.1244286043321187; P=0.1244286043321187; Q=1487.118429386358; R=13384.06586447722; S=0.02858985929294918; T=0.02858985929294918 }
    11 = @{ G=21.811182; H=65.43354600000001; I=0.2297691872894318; J=0.2297691872894318; O=0.345973452289334; P=0.3459734522893341; Q=4134.92942189244; R=37214.36479703196; S=0.07949403895623929; T=0.0794940389562393 }
    12 = @{ G=21.811182; H=65.43354600000001; I=0.2297691872894318; J=0.2297691872894318; M=188.0130056666667; N=564.0390170000001; O=0.3431166302883566; P=0.3431166302883567; Q=4100.785884962699; R=36907.07296466429; S=0.07883762928684415; T=0.07883762928684417 }
    13 = @{ G=21.811182; H=65.43354600000001; I=0.2297691872894318; J=0.2297691872894318; M=102.1836573333333; N=306.550972; O=0.1864813130901906; P=0.1864813130901907; Q=2228.746347522968; R=20058.71712770671; S=0.04284765975339919; T=0.04284765975339919 }
    14 = @{ G=0.3334030000000001; H=1.000209; I=0.003512223058331201; J=0.003512223058331201; M=68.18146900000001; N=204.544407; O=0.1244286043321187; P=0.1244286043321187; Q=22.73190630900701; R=204.587156781063; S=0.0004370210132512368; T=0.0004370210132512367 }
    15 = @{ G=0.3334030000000001; H=1.000209; I=0.003512223058331201; J=0.003512223058331201; O=0.345973452289334; P=0.3459734522893341; Q=63.20601396326; R=568.85412566934; S=0.001215135936701049; T=0.001215135936701049 }
    16 = @{ G=0.3334030000000001; H=1.000209; I=0.003512223058331201; J=0.003512223058331201; M=188.0130056666667; N=564.0390170000001; O=0.3431166302883566; P=0.3431166302883567; Q=62.68410012828368; R=564.1569011545531; S=0.001205102140595668; T=0.001205102140595668 }
    17 = @{ G=0.3334030000000001; H=1.000209; I=0.003512223058331201; J=0.003512223058331201; M=102.1836573333333; N=306.550972; O=0.1864813130901906; P=0.1864813130901907; Q=34.06833790590534; R=306.615041153148; S=0.0006549639677832476; T=0.0006549639677832476 }
}

foreach ($rowNum in $data.Keys) {
    $cols = $data[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
